# Commit: "add relation values to input for SNG from biogas"
#
# Edits the "Unit_relations" sheet:
#  - Methanation (row 6): Input1/Input2 are swapped (h2 / bio_ch4 instead of
#    bio_ch4 / h2), a third input (power) is added, Output2 becomes heat_low
#    (the old co2 / heat_high outputs go away), and relation values are
#    filled in (Relation_In1_In2, Relation_In1_In3, Relation_In_Out,
#    Relation_Out1_Out2).
#  - The CO2_Remover row and the heat_split row are removed entirely, which
#    shifts the water_import and biogas_import rows up two places.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit_relations")

# --- Update the Methanation row (row 6) ---------------------------------
$ws.Range("C6").Value = "h2"
$ws.Range("D6").Value = "bio_ch4"
$ws.Range("E6").Value = "power"
$ws.Range("H6").Value = "heat_low"
$ws.Range("I6").ClearContents()

$ws.Range("K6").Value = 0.86792452830188682
$ws.Range("L6").Value = 46
$ws.Range("N6").Value = 0.5168539325842697
$ws.Range("O6").Value = 8.9

# --- Remove the CO2_Remover row and the heat_split row ------------------
# (water_import / biogas_import shift up from rows 9/10 into rows 7/8)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# --- Restore the view/selection state seen in the saved file ------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("I12").Select()
